$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.873.99'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.815.22'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '308.87'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4666'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3684'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07374'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.79%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8703'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.41'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.793.85'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.25%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.353'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.07059'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.498'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '91.62'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008699'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.73'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '26.921.33'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.343'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.091.08'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.79%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.903'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '150.18'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.171'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.32%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.32'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.88%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.323'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.38%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '115.70'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.34%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08937'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7675'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.10%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.162'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.502'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.901'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.71%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.000'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.087'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.48%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01959'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.20%  '
$ws.Range('E39').Value = '  +1.36%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.937'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.254'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.54%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5313'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.01%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.355'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1662'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.426'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4927'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.12%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.47'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.33%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.000'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '103.69'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.666'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.46%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06285'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.15%  '
